$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.594.81"
$ws.Range("E2").Value = "  -0.33%  "

# Row 3
$ws.Range("D3").Value = "1.841.91"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("Z1").Formula = "=""1.007"""
$ws.Range("Z1").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E4").Value = "  -2.50%  "

# Row 5
$ws.Range("Z1").Formula = "=""316.92"""
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E5").Value = "  -1.61%  "

# Row 6
$ws.Range("Z1").Formula = "=""1.007"""
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E6").Value = "  -2.20%  "

# Row 7
$ws.Range("Z1").Formula = "=""0.4306"""
$ws.Range("Z1").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E7").Value = "  -1.86%  "

# Row 8
$ws.Range("Z1").Formula = "=""0.3726"""
$ws.Range("Z1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E8").Value = "  -1.67%  "

# Row 9
$ws.Range("Z1").Formula = "=""0.07286"""
$ws.Range("Z1").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E9").Value = "  -1.31%  "

# Row 10
$ws.Range("Z1").Formula = "=""0.8708"""
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E10").Value = "  -1.25%  "

# Row 11
$ws.Range("Z1").Formula = "=""21.30"""
$ws.Range("Z1").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E11").Value = "  -1.17%  "

# Row 12
$ws.Range("D12").Value = "1.850.97"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13
$ws.Range("Z1").Formula = "=""6.725"""
$ws.Range("Z1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E13").Value = "  +0.39%  "

# Row 14
$ws.Range("Z1").Formula = "=""5.382"""
$ws.Range("Z1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E14").Value = "  -2.09%  "

# Row 15
$ws.Range("Z1").Formula = "=""0.07091"""
$ws.Range("Z1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E15").Value = "  -0.94%  "

# Row 16
$ws.Range("Z1").Formula = "=""88.59"""
$ws.Range("Z1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E16").Value = "  +4.33%  "

# Row 17
$ws.Range("E17").Value = "  -2.45%  "

# Row 18
$ws.Range("Z1").Formula = "=""0.000008957"""
$ws.Range("Z1").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E18").Value = "  -1.02%  "

# Row 19
$ws.Range("Z1").Formula = "=""1.007"""
$ws.Range("Z1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E19").Value = "  -2.25%  "

# Row 20
$ws.Range("Z1").Formula = "=""15.31"""
$ws.Range("Z1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E20").Value = "  -0.85%  "

# Row 21
$ws.Range("D21").Value = "27.601.65"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("Z1").Formula = "=""5.179"""
$ws.Range("Z1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E22").Value = "  -2.06%  "

# Row 23
$ws.Range("E23").Value = "  -2.77%  "

# Row 24
$ws.Range("D24").Value = "2.074.53"
$ws.Range("E24").Value = "  -0.51%  "

# Row 25
$ws.Range("Z1").Formula = "=""1.964"""
$ws.Range("Z1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E25").Value = "  -4.82%  "

# Row 26
$ws.Range("Z1").Formula = "=""154.28"""
$ws.Range("Z1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E26").Value = "  -2.87%  "

# Row 27
$ws.Range("Z1").Formula = "=""18.52"""
$ws.Range("Z1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E27").Value = "  -0.86%  "

# Row 28
$ws.Range("Z1").Formula = "=""2.159"""
$ws.Range("Z1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E28").Value = "  +8.46%  "

# Row 29
$ws.Range("Z1").Formula = "=""5.311"""
$ws.Range("Z1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("Z1").Formula = "=""117.47"""
$ws.Range("Z1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("Z1").Formula = "=""0.08895"""
$ws.Range("Z1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E31").Value = "  -1.70%  "

# Row 32
$ws.Range("E32").Value = "  +0.55%  "

# Row 33
$ws.Range("Z1").Formula = "=""0.7712"""
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E33").Value = "  +0.12%  "

# Row 34
$ws.Range("Z1").Formula = "=""4.510"""
$ws.Range("Z1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E34").Value = "  -0.95%  "

# Row 35
$ws.Range("Z1").Formula = "=""2.899"""
$ws.Range("Z1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E35").Value = "  -3.45%  "

# Row 36
$ws.Range("Z1").Formula = "=""1.007"""
$ws.Range("Z1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E36").Value = "  -2.35%  "

# Row 37
$ws.Range("E37").Value = "  -2.11%  "

# Row 38
$ws.Range("Z1").Formula = "=""0.01968"""
$ws.Range("Z1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
$ws.Range("Z1").Formula = "=""0.05292"""
$ws.Range("Z1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E39").Value = "  +0.59%  "

# Row 40
$ws.Range("Z1").Formula = "=""2.881"""
$ws.Range("Z1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E40").Value = "  +1.45%  "

# Row 41
$ws.Range("Z1").Formula = "=""7.141"""
$ws.Range("Z1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E41").Value = "  +4.13%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("Z1").Formula = "=""0.5109"""
$ws.Range("Z1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E42").Value = "  -1.21%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("Z1").Formula = "=""0.1683"""
$ws.Range("Z1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E43").Value = "  +0.86%  "

# Row 44
$ws.Range("Z1").Formula = "=""8.744"""
$ws.Range("Z1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E44").Value = "  +0.49%  "

# Row 45
$ws.Range("Z1").Formula = "=""10.65"""
$ws.Range("Z1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("E46").Value = "  -2.98%  "

# Row 47
$ws.Range("Z1").Formula = "=""0.4735"""
$ws.Range("Z1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E47").Value = "  +1.00%  "

# Row 48
$ws.Range("Z1").Formula = "=""0.06441"""
$ws.Range("Z1").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E48").Value = "  -2.30%  "

# Row 49
$ws.Range("E49").Value = "  -2.44%  "

# Row 50
$ws.Range("Z1").Formula = "=""1.676"""
$ws.Range("Z1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E50").Value = "  -1.31%  "

# Row 51
$ws.Range("E51").Value = "  -2.64%  "
